$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at the top of the Ciboulette block (old row 249),
# pushing the existing rows 249:273 down to 251:275.
$ws.Rows("249:250").Insert()

# Fill in the two new rows (249 and 250) with a new weekly data point for
# "Primera" and "Segunda" categories, matching the constant columns used by
# the rest of this dataset block.

# Row 249 - Primera
$ws.Range("A249").Value = 6
$ws.Range("B249").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 44449
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 100112039
$ws.Range("G249").Value = "Ciboulette"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 470
$ws.Range("K249").Value = 1500
$ws.Range("L249").Value = 1500
$ws.Range("M249").Value = 1500
$ws.Range("N249").Value = '$/docena de atados'
$ws.Range("O249").Value = "Región Metropolitana"
$ws.Range("P249").Value = 500
$ws.Range("Q249").Value = 3
$ws.Range("R249").Value = "Hortaliza"

# Row 250 - Segunda
$ws.Range("A250").Value = 6
$ws.Range("B250").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C250").Value = "Metropolitana"
$ws.Range("D250").Value = 44449
$ws.Range("E250").Value = 13
$ws.Range("F250").Value = 100112039
$ws.Range("G250").Value = "Ciboulette"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Segunda"
$ws.Range("J250").Value = 350
$ws.Range("K250").Value = 1000
$ws.Range("L250").Value = 1000
$ws.Range("M250").Value = 1000
$ws.Range("N250").Value = '$/docena de atados'
$ws.Range("O250").Value = "Región Metropolitana"
$ws.Range("P250").Value = 333
$ws.Range("Q250").Value = 3
$ws.Range("R250").Value = "Hortaliza"
